$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 9.77 = 41039.16 pesos`n✅ 41039.16 pesos = 9.7 = 953.67 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update rate values ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 102.39
$ws2.Range("O10").Value = 4202
$ws2.Range("N12").Value = 4231
$ws2.Range("O12").Value = 98.31999999999999
